# Applies the "1401CE19" marksheet update:
#  - recompute the Right/Wrong/Not-Attempt/Max summary numbers (rows 10-12)
#  - fix the "Marking" wrong-answer cell to be numeric instead of text
#  - change the Total/Max cell text from "Absent" to the final score "34/112"
#  - merge the per-question "Student Ans" data (previously spread across the
#    A/D/G column groups) into the single A column, colouring each answer
#    correct/incorrect/unanswered
#  - drop the now-unused 3rd "Student Ans / Correct Ans" column group (G:H)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Summary table (rows 10-12) ----------------------------------------

# Row 10 ("No.") gets the bold/centered title style applied (text unchanged)
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 17
$ws.Range("E10").Value = 28

# Row 11 ("Marking")
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12 ("Total")
$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 36
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "34/112"

# ---- Per-question answers -----------------------------------------------
# Column D of the 2nd group ("Student Ans") picks up the options that used
# to live only in the (now removed) 3rd group, and column A gets the
# consolidated/graded "Student Ans" for every question, styled by
# correctness: correctStyle (green) / incorrectStyle (red).

$ws.Range("D16").Value = "Option A"
$ws.Range("D16").Style = "correctStyle"

$ws.Range("A17").Value = "Option D"
$ws.Range("A17").Style = "correctStyle"
$ws.Range("D17").Value = "Option C"
$ws.Range("D17").Style = "correctStyle"

$ws.Range("A19").Value = "Option C"
$ws.Range("A19").Style = "correctStyle"

$ws.Range("A22").Value = "Option D"
$ws.Range("A22").Style = "correctStyle"

$ws.Range("A25").Value = "Option D"
$ws.Range("A25").Style = "incorrectStyle"

$ws.Range("A26").Value = "Option C"
$ws.Range("A26").Style = "correctStyle"

$ws.Range("A27").Value = "Option A"
$ws.Range("A27").Style = "correctStyle"

$ws.Range("A31").Value = "Option C"
$ws.Range("A31").Style = "incorrectStyle"

$ws.Range("A32").Value = "Option C"
$ws.Range("A32").Style = "correctStyle"

$ws.Range("A36").Value = "Option A"
$ws.Range("A36").Style = "correctStyle"

# Clear out the remaining cells of the no-longer-needed 2nd "Student
# Ans/Correct Ans" rows (they only had a value in D16:E18 which stays).
$ws.Range("D19:E40").Clear()

# ---- Remove the 3rd "Student Ans / Correct Ans" column group (G:H) ------
$ws.Range("G1:H1").EntireColumn.Delete()
